$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "51.774.45"
$ws.Range("E2").Value = "  +4.70%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.754.19"
$ws.Range("E3").Value = "  +4.42%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.13%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.50%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.537"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.15%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.24%  "

$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.95%  "

$ws.Range("B11").Value = "Chainlink"
$ws.Range("C11").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.21"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.64%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.129"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.13%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.40%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.182.71"
$ws.Range("E15").Value = "  +4.04%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.767.91"
$ws.Range("E16").Value = "  +4.69%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.49%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "51.654.40"
$ws.Range("E18").Value = "  +4.48%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.00%  "

$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.25%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +2.21%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.11%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.57%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.57%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.69%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.30%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.140"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.99%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0826"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.00%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.01%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.01"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.61%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.51"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.25%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.57"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.53%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.92%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0345"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.43%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.97%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +14.42%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.108.21"
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.52%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.79%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.65"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.97%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.43%  "
